# Actualizacion automatica del tracker
# Fills in results ("resultado"/"profit") for matches that have finished
# and appends the newly scraped upcoming matches (rows 30-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update resultado/profit columns (G/H) for matches that now have a result
# ---------------------------------------------------------------------------
function Set-Resultado($row, $resultado, $profit) {
    $ws.Cells.Item($row, 7).Value = $resultado
    $ws.Cells.Item($row, 8).Value = $profit
}

Set-Resultado 7  "Acierto" 1.3
Set-Resultado 12 "Acierto" 1.5
Set-Resultado 14 "Fallo"  -1
Set-Resultado 15 "Fallo"  -1
Set-Resultado 16 "Fallo"  -1
Set-Resultado 19 "Acierto" 1.1
Set-Resultado 20 "Fallo"  -1

# ---------------------------------------------------------------------------
# 2) Append the new matches scraped for 2025-09-01 (rows 30-42)
# ---------------------------------------------------------------------------
function Add-Partido($row, $eventId, $fecha, $jugadorA, $jugadorB, $pronostico, $cuota) {
    $ws.Cells.Item($row, 1).Value = $eventId
    # Write the date as a literal text string (not an Excel date serial).
    # Using a formula that evaluates to text and then collapsing it to a
    # plain value via copy/paste keeps the cell a plain text cell without
    # forcing a new "Text" number-format style onto the sheet.
    $ws.Cells.Item($row, 2).Formula = "=""" + $fecha + """"
    $ws.Cells.Item($row, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4163)
    $ws.Cells.Item($row, 3).Value = $jugadorA
    $ws.Cells.Item($row, 4).Value = $jugadorB
    $ws.Cells.Item($row, 5).Value = $pronostico
    $ws.Cells.Item($row, 6).Value = $cuota
}

Add-Partido 30 14580792 "2025-09-01" "Nicolai Budkov Kjaer"      "Borna Ćorić"              "Gana Nicolai Budkov Kjaer"      3.5
Add-Partido 31 14580795 "2025-09-01" "Nikolas Sanchez Izquierdo" "Matej Dodig"              "Gana Matej Dodig"               2.38
Add-Partido 32 14580791 "2025-09-01" "Stefano Travaglia"         "Duje Ajduković"           "Gana Duje Ajduković"            2.5
Add-Partido 33 14580796 "2025-09-01" "Thiago Seyboth Wild"       "Lorenzo Carboni"          "Gana Lorenzo Carboni"           2.63
Add-Partido 34 14581055 "2025-09-01" "Daniel Rincon"             "Pablo Carreño Busta"      "Gana Daniel Rincon"             4
Add-Partido 35 14581058 "2025-09-01" "Dusan Lajovic"             "Elias Ymer"               "Gana Elias Ymer"                3
Add-Partido 36 14581062 "2025-09-01" "Roberto Carballés Baena"   "Genaro Alberto Olivieri"  "Gana Genaro Alberto Olivieri"   3.75
Add-Partido 37 14579767 "2025-09-01" "Daniel Evans"              "Linang Xiao"              "Gana Linang Xiao"               6
Add-Partido 38 14579393 "2025-09-01" "James McCabe"              "Fajing Sun"               "Gana Fajing Sun"                2.75
Add-Partido 39 14579765 "2025-09-01" "Jie Cui"                   "Yu Hsiou Hsu"             "Gana Jie Cui"                   2.75
Add-Partido 40 14579763 "2025-09-01" "Te Rigele"                 "Omar Jasika"              "Gana Te Rigele"                 2.1
Add-Partido 41 14580351 "2025-09-01" "Arthur Fery"               "Michael Geerts"           "Gana Michael Geerts"            3
Add-Partido 42 14580342 "2025-09-01" "Tom Paris"                 "Olle Wallin"              "Gana Olle Wallin"               2.63
